$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$scratch.Formula = '="02-08-2021"'
$scratch.Copy()
$ws.Range("A147").PasteSpecial(-4163)
$ws.Range("B147").Value = 22337

$scratch.Formula = '="03-08-2021"'
$scratch.Copy()
$ws.Range("A148").PasteSpecial(-4163)
$ws.Range("B148").Value = 21351

$scratch.Formula = '="04-08-2021"'
$scratch.Copy()
$ws.Range("A149").PasteSpecial(-4163)
$ws.Range("B149").Value = 22043

$scratch.Formula = '="05-08-2021"'
$scratch.Copy()
$ws.Range("A150").PasteSpecial(-4163)
$ws.Range("B150").Value = 21942

$scratch.Formula = '="06-08-2021"'
$scratch.Copy()
$ws.Range("A151").PasteSpecial(-4163)
$ws.Range("B151").Value = 21122

$scratch.Formula = '="09-08-2021"'
$scratch.Copy()
$ws.Range("A152").PasteSpecial(-4163)
$ws.Range("B152").Value = 26613

$scratch.Formula = '="10-08-2021"'
$scratch.Copy()
$ws.Range("A153").PasteSpecial(-4163)
$ws.Range("B153").Value = 27137

$scratch.Formula = '="11-08-2021"'
$scratch.Copy()
$ws.Range("A154").PasteSpecial(-4163)
$ws.Range("B154").Value = 27260

$scratch.Formula = '="12-08-2021"'
$scratch.Copy()
$ws.Range("A155").PasteSpecial(-4163)
$ws.Range("B155").Value = 27140

$scratch.Formula = '="13-08-2021"'
$scratch.Copy()
$ws.Range("A156").PasteSpecial(-4163)
$ws.Range("B156").Value = 26028

$scratch.Formula = '="16-08-2021"'
$scratch.Copy()
$ws.Range("A157").PasteSpecial(-4163)
$ws.Range("B157").Value = 26318

$scratch.Formula = '="17-08-2021"'
$scratch.Copy()
$ws.Range("A158").PasteSpecial(-4163)
$ws.Range("B158").Value = 26351

$scratch.Formula = '="18-08-2021"'
$scratch.Copy()
$ws.Range("A159").PasteSpecial(-4163)
$ws.Range("B159").Value = 26435

$scratch.Formula = '="19-08-2021"'
$scratch.Copy()
$ws.Range("A160").PasteSpecial(-4163)
$ws.Range("B160").Value = 25174

$scratch.Formula = '="20-08-2021"'
$scratch.Copy()
$ws.Range("A161").PasteSpecial(-4163)
$ws.Range("B161").Value = 25592

$scratch.Formula = '="23-08-2021"'
$scratch.Copy()
$ws.Range("A162").PasteSpecial(-4163)
$ws.Range("B162").Value = 24987

$scratch.Formula = '="24-08-2021"'
$scratch.Copy()
$ws.Range("A163").PasteSpecial(-4163)
$ws.Range("B163").Value = 24934

$scratch.Formula = '="25-08-2021"'
$scratch.Copy()
$ws.Range("A164").PasteSpecial(-4163)
$ws.Range("B164").Value = 25574

$scratch.Formula = '="26-08-2021"'
$scratch.Copy()
$ws.Range("A165").PasteSpecial(-4163)
$ws.Range("B165").Value = 25498

$scratch.Formula = '="27-08-2021"'
$scratch.Copy()
$ws.Range("A166").PasteSpecial(-4163)
$ws.Range("B166").Value = 25248

$scratch.Formula = '="30-08-2021"'
$scratch.Copy()
$ws.Range("A167").PasteSpecial(-4163)
$ws.Range("B167").Value = 26492

$scratch.Formula = '="31-08-2021"'
$scratch.Copy()
$ws.Range("A168").PasteSpecial(-4163)
$ws.Range("B168").Value = 28297

$scratch.ClearContents()
"done"